$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1) to short machine-friendly names
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'
$ws.Range('B14').Value = 'Amatenango De La Frontera'
$ws.Range('B19').Value = 'Comitán De Domínguez'
$ws.Range('B31').Value = 'Mazapa De Madero'
$ws.Range('B41').Value = 'San Cristóbal De Las Casas'
$ws.Range('A58').Value = 'Ciudad De México'
$ws.Range('A78').Value = 'Estado De México'
$ws.Range('B78').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B80').Value = 'Atizapán De Zaragoza'
$ws.Range('B85').Value = 'Ecatepec De Morelos'
$ws.Range('B90').Value = 'Naucalpan De Juárez'
$ws.Range('B92').Value = 'San Felipe Del Progreso'
$ws.Range('B97').Value = 'Tlalnepantla De Baz'
$ws.Range('B100').Value = 'Villa De Allende'
$ws.Range('B104').Value = 'Apaseo El Alto'
$ws.Range('B107').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B113').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B114').Value = 'Valle De Santiago'
$ws.Range('B117').Value = 'Acapulco De Juárez'
$ws.Range('B120').Value = 'Ayutla De Los Libres'
$ws.Range('B122').Value = 'Chilapa De Álvarez'
$ws.Range('B124').Value = 'Coyuca De Benítez'
$ws.Range('B125').Value = 'Coyuca De Catalán'
$ws.Range('B132').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B134').Value = 'Tlapa De Comonfort'
$ws.Range('B149').Value = 'Tenango De Doria'
$ws.Range('B150').Value = 'Tula De Allende'
$ws.Range('B151').Value = 'Tulancingo De Bravo'
$ws.Range('B155').Value = 'Ahualulco De Mercado'
$ws.Range('B157').Value = 'Autlán De Navarro'
$ws.Range('B162').Value = 'Cuautitlán De García Barragán'
$ws.Range('B168').Value = 'Tizapán El Alto'
$ws.Range('B169').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B172').Value = 'Unión De San Antonio'
$ws.Range('B173').Value = 'Unión De Tula'
$ws.Range('B175').Value = 'Zapotlán El Grande'
$ws.Range('B198').Value = 'Tetela Del Volcán'
$ws.Range('B206').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B209').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B210').Value = 'Ixtlán De Juárez'
$ws.Range('B211').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B213').Value = 'Oaxaca De Juárez'
$ws.Range('B214').Value = 'Pinotepa De Don Luis'
$ws.Range('B215').Value = 'Putla Villa De Guerrero'
$ws.Range('B220').Value = 'San Dionisio Del Mar'
$ws.Range('B231').Value = 'Totontepec Villa De Morelos'
$ws.Range('B244').Value = 'Izúcar De Matamoros'
$ws.Range('B251').Value = 'San Salvador El Verde'
$ws.Range('B253').Value = 'Tepanco De López'
$ws.Range('B263').Value = 'Xayacatlán De Bravo'
$ws.Range('B269').Value = 'Landa De Matamoros'
$ws.Range('B270').Value = 'Pinal De Amoles'
$ws.Range('B272').Value = 'San Juan Del Río'
$ws.Range('B283').Value = 'Villa De Reyes'
$ws.Range('B311').Value = 'Tepetitla De Lardizábal'
$ws.Range('B317').Value = 'Castillo De Teayo'
$ws.Range('B325').Value = 'Ignacio De La Llave'
$ws.Range('B328').Value = 'Juchique De Ferrer'
$ws.Range('B332').Value = 'Ozuluama De Mascareñas'
$ws.Range('B334').Value = 'Sayula De Alemán'

# Remove the trailing footer/metadata rows (355-359), shrinking the used range to A1:D353
$ws.Range("A355:A359").EntireRow.Delete() | Out-Null
